$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.88000000000014
$ws.Range("H2").Value = [double]"2.413528314402514e-16"
$ws.Range("K2").Value = 39.6738796845174
$ws.Range("L2").Value = "[31.33223298719693, 48.01552638183786]"
$ws.Range("O2").Value = 1.754763464167272
$ws.Range("P2").Value = "[1.5157634224527339, 1.9937635058818106]"
$ws.Range("S2").Value = 57.92916248496812
$ws.Range("T2").Value = "[52.752583525290774, 63.105741444645474]"
$ws.Range("W2").Value = 16.49009009009019
$ws.Range("X2").Value = 15.61977977977987
$ws.Range("Y2").Value = 17.3604004004005

# Row 3 updates
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 25.8500000000006
$ws.Range("H3").Value = [double]"2.413528314402514e-16"
$ws.Range("K3").Value = 46.72982678577304
$ws.Range("L3").Value = "[38.259344185007414, 55.20030938653867]"
$ws.Range("O3").Value = 0.1698158191129613
$ws.Range("P3").Value = "[-0.031447373909808896, 0.37107901213573147]"
$ws.Range("Q3").Value = 0.09785889909533241
$ws.Range("R3").Value = 0.09785889909533241
$ws.Range("S3").Value = 61.30705446738548
$ws.Range("T3").Value = "[56.01242259532043, 66.60168633945054]"
$ws.Range("W3").Value = 25.15135135135194
$ws.Range("X3").Value = 24.32332332332389
$ws.Range("Y3").Value = 25.97937937937999
